# Refresh the crypto price table (coinranking.com snapshot).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# Two coin pairs swapped ranking position between snapshots
# (rows 15/16 and rows 36/37/40/42), so B/C/D/E are all rewritten
# for those rows; everywhere else only D and/or E changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ Column letter = new value }. Column D sometimes holds
# purely-numeric-looking text (e.g. '1.00', '0.0460', '3.595.21');
# writing it with Range.Value would let Excel coerce it to a
# number and silently drop the formatting, so those cells are
# pre-formatted as Text before the value is written.
$updates = [ordered]@{
    2 = @{ D='70.137.32'; E='  +0.50%  ' }
    3 = @{ D='3.609.53'; E='  +3.32%  ' }
    4 = @{ E='  -0.16%  ' }
    5 = @{ D='604.95'; E='  +0.07%  ' }
    6 = @{ D='195.81'; E='  -1.40%  ' }
    7 = @{ D='0.626'; E='  +0.02%  ' }
    8 = @{ E='  +0.05%  ' }
    9 = @{ E='  -2.05%  ' }
    10 = @{ E='  -0.43%  ' }
    11 = @{ D='54.02'; E='  -0.43%  ' }
    12 = @{ E='  +0.33%  ' }
    13 = @{ D='9.56'; E='  -0.10%  ' }
    14 = @{ D='4.174.71'; E='  +2.91%  ' }
    15 = @{ B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='13.09'; E='  +3.64%  ' }
    16 = @{ B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='596.69'; E='  +0.12%  ' }
    17 = @{ D='19.33'; E='  +1.97%  ' }
    18 = @{ D='70.349.67'; E='  +0.69%  ' }
    19 = @{ D='3.595.21'; E='  +3.30%  ' }
    20 = @{ D='0.122'; E='  +1.67%  ' }
    21 = @{ D='0.998'; E='  +1.09%  ' }
    22 = @{ D='17.91'; E='  -0.03%  ' }
    23 = @{ E='  +2.12%  ' }
    24 = @{ D='102.84'; E='  -0.61%  ' }
    25 = @{ E='  +0.14%  ' }
    26 = @{ D='3.09'; E='  -0.54%  ' }
    27 = @{ E='  -0.54%  ' }
    28 = @{ E='  -2.31%  ' }
    29 = @{ D='34.16'; E='  +1.84%  ' }
    30 = @{ E='  -1.49%  ' }
    31 = @{ E='  -4.90%  ' }
    32 = @{ D='12.35'; E='  -3.31%  ' }
    33 = @{ D='0.117'; E='  +0.86%  ' }
    34 = @{ D='63.58'; E='  -0.16%  ' }
    35 = @{ D='3.903.00'; E='  +5.31%  ' }
    36 = @{ B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.0₃0851'; E='  +5.23%  ' }
    37 = @{ B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='3.21'; E='  +7.97%  ' }
    38 = @{ D='530.64'; E='  +1.81%  ' }
    39 = @{ E='  +0.08%  ' }
    40 = @{ B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.395'; E='  +1.00%  ' }
    41 = @{ D='37.13'; E='  +0.93%  ' }
    42 = @{ B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='3.59'; E='  +1.05%  ' }
    43 = @{ E='  -2.32%  ' }
    44 = @{ D='0.0460'; E='  +0.36%  ' }
    45 = @{ D='2.87'; E='  +0.56%  ' }
    46 = @{ D='3.36'; E='  +1.55%  ' }
    47 = @{ D='0.141'; E='  +0.72%  ' }
    48 = @{ D='8.62'; E='  -1.56%  ' }
    49 = @{ D='1.00'; E='  -0.08%  ' }
    50 = @{ E='  +3.44%  ' }
    51 = @{ D='1.30' }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cell = $ws.Range("$col$row")
        $newValue = $updates[$row][$col]
        if ($col -eq 'D') {
            $cell.NumberFormat = '@'
        }
        $cell.Value = $newValue
    }
}
